# Condition1.xlsx edit: switch path separators from "/" to "\" in the
# SoundName/CSName/UCSName path columns, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace forward slashes with backslashes in the used data range (A2:C31
# holds the "Condition/xxx" and "Sound/xxx" style path strings).
$used = $ws.UsedRange
$used.Replace("/", "\")

# Move the active selection to A28 (matches the saved sheetView state).
$ws.Range("A28").Select()
